$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.475.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.77%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.838.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.31%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.86%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.71%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4611"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3826"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.42"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07909"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9764"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.830.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.909"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.055"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06640"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001029"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.467.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.372"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.302"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.071"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.338"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9568"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09300"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.567"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.251"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.322"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05952"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02196"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.086"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.161"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5813"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1844"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.250"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5500"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.99%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.874"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06654"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.044"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.62%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.003"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.74%  "

$ws.Range("E51").Value = "  -1.54%  "

